# Add a new worksheet "attSearchResultPagePromoMessage" after "phonesAndDevices",
# populate it with the two new assertion strings, make it the active/selected
# sheet with H32 selected, and remove the tabSelected flag from sheet1.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "attSearchResultPagePromoMessage"

$ws2.Range("A1").Value = "att search results test assert"
$ws2.Range("A2").Value = "Shop the latest offers designed for individual first responders and those that support them. See how FirstNet can help save you money."

$ws2.Range("H32").Select() | Out-Null
